$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (E1) onto the new
# header cell F1 so it reuses the same style as the rest of row 1, then set
# its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Add the new data cell with the model description (default/no special style).
$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"
